# Apply updated currentAveragePrice / LevePrice / LeveProfit values
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets (scheduled price refresh).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 345.69232
$ws.Range("I96").Value = 364.5
$ws.Range("J96").Value = 283
$ws.Range("K96").Value = 1093.5
$ws.Range("L96").Value = 849
$ws.Range("M96").Value = 279.5
$ws.Range("N96").Value = -3595
$ws.Range("H98").Value = 1923.625
$ws.Range("I98").Value = 1842.7142
$ws.Range("J98").Value = 2490
$ws.Range("K98").Value = 1842.7142
$ws.Range("L98").Value = 2490
$ws.Range("M98").Value = -344.7141999999999
$ws.Range("N98").Value = -5486
$ws.Range("H100").Value = 2105.889
$ws.Range("I100").Value = 2396
$ws.Range("K100").Value = 2396
$ws.Range("M100").Value = -1855
$ws.Range("H106").Value = 1511
$ws.Range("I106").Value = 1388.75
$ws.Range("K106").Value = 1388.75
$ws.Range("M106").Value = -757.75
$ws.Range("H122").Value = 1923.625
$ws.Range("I122").Value = 1842.7142
$ws.Range("J122").Value = 2490
$ws.Range("K122").Value = 5528.142599999999
$ws.Range("L122").Value = 7470
$ws.Range("M122").Value = -3078.142599999999
$ws.Range("N122").Value = -12370
$ws.Range("H125").Value = 254974.5
$ws.Range("I125").Value = 4000
$ws.Range("K125").Value = 36000
$ws.Range("M125").Value = -33540

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 579.25
$ws.Range("I2").Value = 489.18182
$ws.Range("J2").Value = 1570
$ws.Range("K2").Value = 489.18182
$ws.Range("L2").Value = 1570
$ws.Range("M2").Value = -376.18182
$ws.Range("N2").Value = -1796
$ws.Range("H32").Value = 2181.24
$ws.Range("I32").Value = 1855.4584
$ws.Range("K32").Value = 1855.4584
$ws.Range("M32").Value = -1568.4584
$ws.Range("H63").Value = 1500
$ws.Range("I63").Value = 1500
$ws.Range("K63").Value = 1500
$ws.Range("M63").Value = -814
$ws.Range("H66").Value = 1500
$ws.Range("I66").Value = 1500
$ws.Range("K66").Value = 7500
$ws.Range("M66").Value = -4068
$ws.Range("H74").Value = 676.25
$ws.Range("I74").Value = 572.8570999999999
$ws.Range("K74").Value = 572.8570999999999
$ws.Range("M74").Value = 301.1429000000001
$ws.Range("H77").Value = 676.25
$ws.Range("I77").Value = 572.8570999999999
$ws.Range("K77").Value = 2864.2855
$ws.Range("M77").Value = 1503.7145
$ws.Range("H116").Value = 579.25
$ws.Range("I116").Value = 489.18182
$ws.Range("J116").Value = 1570
$ws.Range("K116").Value = 489.18182
$ws.Range("L116").Value = 1570
$ws.Range("M116").Value = 1804.81818
$ws.Range("N116").Value = -6158
$ws.Range("H122").Value = 2424.3333
$ws.Range("I122").Value = 1320.5
$ws.Range("K122").Value = 3961.5
$ws.Range("M122").Value = -1511.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 579.25
$ws.Range("I3").Value = 489.18182
$ws.Range("J3").Value = 1570
$ws.Range("K3").Value = 489.18182
$ws.Range("L3").Value = 1570
$ws.Range("M3").Value = -375.18182
$ws.Range("N3").Value = -1798

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 47.555557
$ws.Range("I7").Value = 41
$ws.Range("K7").Value = 41
$ws.Range("M7").Value = 72
$ws.Range("H16").Value = 1364.2858
$ws.Range("I16").Value = 1070
$ws.Range("K16").Value = 1070
$ws.Range("M16").Value = -783
$ws.Range("H22").Value = 613.4286
$ws.Range("J22").Value = 735.75
$ws.Range("L22").Value = 735.75
$ws.Range("N22").Value = -1435.75
$ws.Range("H41").Value = 5029.5
$ws.Range("H99").Value = 4832.6665
$ws.Range("I99").Value = 4148.5
$ws.Range("K99").Value = 4148.5
$ws.Range("M99").Value = -2650.5
$ws.Range("H105").Value = 1098.1818
$ws.Range("I105").Value = 910
$ws.Range("J105").Value = 1945
$ws.Range("K105").Value = 910
$ws.Range("L105").Value = 1945
$ws.Range("M105").Value = 837
$ws.Range("N105").Value = -5439
$ws.Range("H107").Value = 454
$ws.Range("I107").Value = 434.9
$ws.Range("J107").Value = 645
$ws.Range("K107").Value = 434.9
$ws.Range("L107").Value = 645
$ws.Range("M107").Value = 1485.1
$ws.Range("N107").Value = -4485
$ws.Range("H113").Value = 1364.2858
$ws.Range("I113").Value = 1070
$ws.Range("K113").Value = 1070
$ws.Range("M113").Value = 1100
$ws.Range("H122").Value = 1097.4
$ws.Range("I122").Value = 1269.8572
$ws.Range("J122").Value = 946.5
$ws.Range("K122").Value = 3809.5716
$ws.Range("L122").Value = 2839.5
$ws.Range("M122").Value = -1359.5716
$ws.Range("N122").Value = -7739.5
$ws.Range("H126").Value = 4832.6665
$ws.Range("I126").Value = 4148.5
$ws.Range("K126").Value = 12445.5
$ws.Range("M126").Value = -9975.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2406.5334
$ws.Range("I4").Value = 2430
$ws.Range("J4").Value = 2359.6
$ws.Range("K4").Value = 7290
$ws.Range("L4").Value = 7078.799999999999
$ws.Range("M4").Value = -7178
$ws.Range("N4").Value = -7302.799999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 53250
$ws.Range("J26").Value = 69875
$ws.Range("L26").Value = 69875
$ws.Range("N26").Value = -70435
$ws.Range("H50").Value = 53250
$ws.Range("J50").Value = 69875
$ws.Range("L50").Value = 69875
$ws.Range("N50").Value = -70871
$ws.Range("H122").Value = 1468.1666
$ws.Range("I122").Value = 1444.8
$ws.Range("K122").Value = 4334.4
$ws.Range("M122").Value = -1884.4
$ws.Range("H132").Value = 2539.75
$ws.Range("I132").Value = 2185
$ws.Range("K132").Value = 6555
$ws.Range("M132").Value = -4025

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 21761.084
$ws.Range("I7").Value = 19329.75
$ws.Range("K7").Value = 19329.75
$ws.Range("M7").Value = -19217.75
$ws.Range("H126").Value = 21761.084
$ws.Range("I126").Value = 19329.75
$ws.Range("K126").Value = 57989.25
$ws.Range("M126").Value = -55519.25
$ws.Range("H132").Value = 7291.5
$ws.Range("I132").Value = 9657.333000000001
$ws.Range("K132").Value = 28971.999
$ws.Range("M132").Value = -26441.999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 2778.8125
$ws.Range("J107").Value = 1997.5714
$ws.Range("L107").Value = 5992.7142
$ws.Range("N107").Value = -9832.7142
$ws.Range("H113").Value = 279.16666
$ws.Range("I113").Value = 242.25
$ws.Range("K113").Value = 726.75
$ws.Range("M113").Value = 1443.25
